$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that needs to move from
# 2023-09-10 (45179) to 2023-09-11 (45180) for every data row (2..367).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 367 }

$ws.Range("C2:C$lastRow").Value = 45180
